$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dimmer")

# Cell D3: import list - rename TimedDimmer to DimmerCommand
$ws.Range("D3").Value = "org.eclipse.smarthome.core.items.Item,`norg.eclipse.smarthome.core.library.types.IncreaseDecreaseType,`norg.eclipse.smarthome.core.library.types.OpenClosedType,`ncom.incquerylabs.smarthome.eventbus.api.events.ItemStateChangedEvent,`ncom.incquerylabs.smarthome.eventbus.api.IEventPublisher,`ncom.incquerylabs.smarthome.eventbus.ruleengine.droolshomeio.DimmerCommand"

# Cell F12: action on start rule - rename timedCommand/TimedDimmer to startComplexCommand/DimmerCommand
$ws.Range("F12").Value = "openhab.startComplexCommand(DimmerCommand.create(`$light, `$param, 2 , 1000, 50));"

# Cells F28, F29, F30: action on stop rule - rename stopTimedCommand to stopComplexCommand
$ws.Range("F28").Value = "openhab.stopComplexCommand(`$light);"
$ws.Range("F29").Value = "openhab.stopComplexCommand(`$light);"
$ws.Range("F30").Value = "openhab.stopComplexCommand(`$light);"
